# edit.ps1
#
# Applies two text edits to Business_Case.docx:
#
# 1. "Snacks aus dem Automaten werden..." paragraph: the sentence was
#    previously split across two runs around a (hidden) "_GoBack"
#    bookmark. The bookmark is removed and the two runs collapse back
#    into a single run holding the whole sentence.
#
# 2. "Es ist nicht mit mehr Wartungsaufwand..." paragraph: "einem
#    personalbetrieben Kino" becomes "einem personalbetriebenen Kino"
#    (an "en" is typed in after "personalbetrieben") and Word's
#    "_GoBack" bookmark (marking the most recent edit point) is now
#    left right after the newly typed "en", before " Kino.".  The
#    surrounding grammar-check proofErr markers are gone too, since
#    that run was touched by the edit.

$d = $word.ActiveDocument

function Get-ParagraphContaining([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Snacks paragraph: merge the two runs back together and drop the
#    bookmark that used to sit between them.
# ---------------------------------------------------------------------

$snacksOld = "Snacks aus dem Automaten werden mithilfe der App bestellt, bereitgestellt und abgerechnet."
$snacksPara = Get-ParagraphContaining "Snacks aus dem Automaten"
$snacksRange = $snacksPara.Range
$snacksRange.Find.Execute($snacksOld, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $snacksOld, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Wartungsaufwand paragraph: insert "en" after "personalbetrieben"
#    and re-home the "_GoBack" bookmark right after it.
# ---------------------------------------------------------------------

$wartOld = "Es ist nicht mit mehr Wartungsaufwand zu rechnen, als bei einem personalbetrieben Kino."
$wartNew = "Es ist nicht mit mehr Wartungsaufwand zu rechnen, als bei einem personalbetriebenen Kino."
$wartPara = Get-ParagraphContaining "Wartungsaufwand"
$wartRange = $wartPara.Range
$wartRange.Find.Execute($wartOld, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $wartNew, 2) | Out-Null

# Re-fetch the paragraph/range after the replace above.
$wartPara = Get-ParagraphContaining "Wartungsaufwand"
$paraStart = $wartPara.Range.Start
$stem = "Es ist nicht mit mehr Wartungsaufwand zu rechnen, als bei einem personalbetrieben"

# Split point between "...personalbetrieben" and the newly typed "en".
$splitPos = $paraStart + $stem.Length
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("SplitMarker", $splitRange) | Out-Null

# Split point (and bookmark home) after "...personalbetriebenen", i.e.
# right before " Kino.".
$bookmarkPos = $paraStart + $stem.Length + 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# The temporary marker only existed to force the run split between
# "personalbetrieben" and "en"; Word's real split doesn't keep a
# bookmark there, so remove it again.
$d.Bookmarks.Item("SplitMarker").Delete()

Write-Output "Snacks paragraph: $($snacksPara.Range.Text)"
Write-Output "Wartungsaufwand paragraph: $($(Get-ParagraphContaining "Wartungsaufwand").Range.Text)"
